$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D prices keep their exact text representation (e.g. "1.003",
# "0.5800", "28.057.03") instead of being auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '27.992.22'
$ws.Range("E2").Value = '  -3.80%  '
$ws.Range("D3").Value = '1.744.90'
$ws.Range("E3").Value = '  -4.29%  '
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").Value = '225.73'
$ws.Range("E5").Value = '  -3.85%  '
$ws.Range("D6").Value = '0.5800'
$ws.Range("E6").Value = '  -3.20%  '
$ws.Range("D7").Value = '1.004'
$ws.Range("E7").Value = '  +0.15%  '
$ws.Range("D8").Value = '0.2703'
$ws.Range("E8").Value = '  -1.74%  '
$ws.Range("D9").Value = '23.06'
$ws.Range("E9").Value = '  -0.81%  '
$ws.Range("D10").Value = '0.06558'
$ws.Range("E10").Value = '  -5.62%  '
$ws.Range("D11").Value = '0.07507'
$ws.Range("E11").Value = '  -1.06%  '
$ws.Range("D12").Value = '1.748.69'
$ws.Range("E12").Value = '  -4.35%  '
$ws.Range("D13").Value = '4.701'
$ws.Range("E13").Value = '  -0.75%  '
$ws.Range("D14").Value = '0.6019'
$ws.Range("E14").Value = '  -3.82%  '
$ws.Range("D15").Value = '1.983.22'
$ws.Range("E15").Value = '  -4.22%  '
$ws.Range("D16").Value = '73.80'
$ws.Range("E16").Value = '  -4.49%  '
$ws.Range("D17").Value = '0.000008622'
$ws.Range("E17").Value = '  -12.09%  '
$ws.Range("D18").Value = '28.003.64'
$ws.Range("E18").Value = '  -2.73%  '
$ws.Range("D19").Value = '5.311'
$ws.Range("E19").Value = '  -4.87%  '
$ws.Range("D20").Value = '1.003'
$ws.Range("E20").Value = '  -0.06%  '
$ws.Range("D21").Value = '203.77'
$ws.Range("E21").Value = '  -5.77%  '
$ws.Range("D22").Value = '11.24'
$ws.Range("E22").Value = '  -2.50%  '
$ws.Range("D23").Value = '6.629'
$ws.Range("E23").Value = '  -3.52%  '
$ws.Range("D24").Value = '1.004'
$ws.Range("E24").Value = '  +0.13%  '
$ws.Range("D25").Value = '150.37'
$ws.Range("E25").Value = '  -3.85%  '
$ws.Range("D26").Value = '8.038'
$ws.Range("E26").Value = '  +1.42%  '
$ws.Range("D27").Value = '0.1230'
$ws.Range("E27").Value = '  -4.35%  '
$ws.Range("D28").Value = '16.05'
$ws.Range("E28").Value = '  -2.41%  '
$ws.Range("E29").Value = '  -2.83%  '
$ws.Range("D30").Value = '0.06074'
$ws.Range("E30").Value = '  -5.05%  '
$ws.Range("D31").Value = '1.385'
$ws.Range("E31").Value = '  -3.78%  '
$ws.Range("E32").Value = '  -2.62%  '
$ws.Range("D33").Value = '3.699'
$ws.Range("E33").Value = '  -1.42%  '
$ws.Range("D34").Value = '1.669'
$ws.Range("E34").Value = '  -3.31%  '
$ws.Range("D35").Value = '1.032'
$ws.Range("E35").Value = '  -5.22%  '
$ws.Range("D36").Value = '0.6329'
$ws.Range("E36").Value = '  -2.14%  '
$ws.Range("D37").Value = '2.452'
$ws.Range("E37").Value = '  -3.18%  '
$ws.Range("D38").Value = '2.712'
$ws.Range("E38").Value = '  -0.84%  '
$ws.Range("D39").Value = '0.01669'
$ws.Range("E39").Value = '  -4.34%  '
$ws.Range("D40").Value = '6.276'
$ws.Range("E40").Value = '  -3.94%  '
$ws.Range("D41").Value = '1.123.41'
$ws.Range("E41").Value = '  -1.86%  '
$ws.Range("D42").Value = '0.8638'
$ws.Range("E42").Value = '  -2.42%  '
$ws.Range("D43").Value = '1.004'
$ws.Range("E43").Value = '  +0.22%  '
$ws.Range("D44").Value = '99.29'
$ws.Range("E44").Value = '  -1.04%  '
$ws.Range("D45").Value = '1.897.92'
$ws.Range("E45").Value = '  -4.17%  '
$ws.Range("D46").Value = '59.05'
$ws.Range("E46").Value = '  -4.07%  '
$ws.Range("D47").Value = '1.572'
$ws.Range("E47").Value = '  -2.04%  '
$ws.Range("E48").Value = '  -5.00%  '
$ws.Range("D49").Value = '8.235'
$ws.Range("E49").Value = '  -2.94%  '
$ws.Range("D50").Value = '0.05391'
$ws.Range("D51").Value = '0.4430'
$ws.Range("E51").Value = '  -2.24%  '
